$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Fill in the new "day 8" entry (row 9)
# Order matches how the new shared strings were appended (image_url, title, text)
$ws.Range("E9").Value = "https://www.visualstatements.net/wp-content/uploads/2023/02/Spruch3-Von-nun-an-egal-was-andere-denken.png"
$ws.Range("B9").Value = "Ist egal was andere denken"
$ws.Range("C9").Value = "Heute habe ich mich mit Corina Brüngger in St.Gallen getroffen und wir haben uns über unsere Visionen unterhalten und gechallenged. Es war so ein schöner Austausch und wir beide konnten uns gegenseitig weiterhelfen. Wir haben uns mega offen unterhalten und es war uns egal, was die rundherum gedacht haben. Ich hatte das Gefühl, dass einige uns beobachtet und zugehört haben, aber es war mir egal. Es gab auch schon Zeiten, da wäre mir das etwas unangenehm gewesen. Was würden wohl die anderen denken? Gottseidank ist das ab jetzt anders."
$ws.Range("D9").Value = "Mut"

# Move the active selection to D10, as it ended up after the edit
$ws.Range("D10").Select()
